# Overwrite player settings from command line: extend the sweep table with
# two more parameter points (0.25, 0.1) and a new "find_root_node" metric
# row, then (re)plot the three series on a line chart anchored to the right
# of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 6: input variable (pts ratio category axis) -----------------------
$ws.Range("B6").Value = 2.5
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 1.5
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.25
$ws.Range("H6").Value = 0.1

# --- Row 7: pts ratio --------------------------------------------------
$ws.Range("B7").Value = 0.466354
$ws.Range("C7").Value = 0.484517
$ws.Range("D7").Value = 0.476697
$ws.Range("E7").Value = 0.474178
$ws.Range("F7").Value = 0.471592
$ws.Range("G7").Value = 0.479281
$ws.Range("H7").Value = 0.437647

# --- Row 8: found terminal node (formulas) ------------------------------
$ws.Range("B8").Formula = "=39104/(52219+39104)"
$ws.Range("C8").Formula = "=38288/(54901+38288)"
$ws.Range("D8").Formula = "=39068/(54982+39068)"
$ws.Range("E8").Formula = "=34516/(58410+34516)"
$ws.Range("F8").Formula = "=33823/(55625+33823)"
$ws.Range("G8").Formula = "=28047/(61526+28047)"
$ws.Range("H8").Formula = "=17036/(77675+17036)"

# --- Row 9 (new): found root node (formulas) ----------------------------
$ws.Range("A9").Value = "find_root_node"
$ws.Range("B9").Formula = "=900/(180+900)"
$ws.Range("C9").Formula = "=925/(190+925)"
$ws.Range("D9").Formula = "=910/(172+910)"
$ws.Range("E9").Formula = "=944/(234+944)"
$ws.Range("F9").Formula = "=785/(411+785)"
$ws.Range("G9").Formula = "=465/(609+465)"
$ws.Range("H9").Formula = "=454/(754+454)"

# --- Chart: line chart of the three metric rows against row 6 --------------
$co = $ws.ChartObjects().Add(370000, 0, 2350000, 3500000)
$chart = $co.Chart
$chart.ChartType = 4

$s1 = $chart.SeriesCollection().NewSeries()
$s1.Name = "pts ratio"
$s1.XValues = $ws.Range("B6:H6")
$s1.Values = $ws.Range("B7:H7")

$s2 = $chart.SeriesCollection().NewSeries()
$s2.Name = "found terminal node"
$s2.XValues = $ws.Range("B6:H6")
$s2.Values = $ws.Range("B8:H8")

$s3 = $chart.SeriesCollection().NewSeries()
$s3.Name = "found root node"
$s3.XValues = $ws.Range("B6:H6")
$s3.Values = $ws.Range("B9:H9")

$chart.HasLegend = $true
$chart.Legend.Position = -4152

# --- Restore the selection that the original author left on the sheet ------
$ws.Range("I26").Select()
